# Astrolabe measurements: replace the "Date" column's date-serial values
# with plain day-of-month numbers (and drop the date number-format from
# those cells), add a missing Altitude 3 reading for row 23, widen the two
# new columns that reading creates visual room for, and leave the
# selection on the cell that was being edited (F24).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column C ("Date"): every data row 2-51 keeps/gets a plain day number
# (16, 17, 22 or 29) instead of a full date serial, and loses the custom
# date number format it carried (style index 1 -> default/Normal).
$dayForRow = @{}
for ($r = 2;  $r -le 9;  $r++) { $dayForRow[$r] = 16 }
$dayForRow[10] = 16
$dayForRow[11] = 16
$dayForRow[12] = 16
$dayForRow[13] = 16
for ($r = 14; $r -le 31; $r++) { $dayForRow[$r] = 17 }
for ($r = 32; $r -le 44; $r++) { $dayForRow[$r] = 22 }
for ($r = 45; $r -le 51; $r++) { $dayForRow[$r] = 29 }

for ($r = 2; $r -le 51; $r++) {
    $ws.Cells.Item($r, 3).Value = $dayForRow[$r]
}

# Strip the date number format from the whole column range in one shot so
# the workbook's style table stops referencing it.
$ws.Range("C2:C51").Style = "Normal"

# --- Row 23 was missing its third altitude reading; fill it in. The
# dependent Average/StdDev formulas (shared formulas in G/H) recalculate
# automatically.
$ws.Range("F23").Value = 48

# --- New columns F and G now hold real data, so give them explicit
# widths (Excel OM widths are in "characters"; the underlying engine
# quantizes to 1/6-character steps, so these inputs were chosen to land
# on/near the target stored widths of 17 and 17.7109375 respectively).
$ws.Columns.Item(6).ColumnWidth = 16.17
$ws.Columns.Item(7).ColumnWidth = 16.88

# --- Leave the selection where the edit happened.
$ws.Range("F24").Select()

# Best-effort: also nudge the window scroll position toward row 5 so the
# view roughly matches what was scrolled into frame while editing.
$win = $excel.ActiveWindow
$win.ScrollRow = 5
$win.ScrollColumn = 1
